# B6-PowerPoint.pptx edit
#   1. Re-style the three tables (slides 14-16) from the custom
#      "Table_0" style {D68BB783-94E6-4599-9288-F0F85294761D} to
#      {242320B5-FCB1-4832-8159-70C79D9FADB5}.
#   2. Swap the deck's two embedded themes ("Office Theme" <-> "Integral")
#      so the design applied to the slide master picks up the colours
#      that used to live on the notes-master-only theme.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# 1) Table styles
# ---------------------------------------------------------------
$newTableStyleId = "{242320B5-FCB1-4832-8159-70C79D9FADB5}"

foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyleId)
        }
    }
}

# ---------------------------------------------------------------
# 2) Theme colour swap (Office Theme <-> Integral / "Red Violet")
# ---------------------------------------------------------------
function Get-RGBValue($r, $g, $b) {
    return $b * 65536 + $g * 256 + $r
}

# The colours the "Office Theme" clrScheme used before the edit - these
# become the live design's colours after the swap.
$officeColors = @(
    (Get-RGBValue 0x00 0x00 0x00),  # dk1
    (Get-RGBValue 0xFF 0xFF 0xFF),  # lt1
    (Get-RGBValue 0x44 0x54 0x6A),  # dk2
    (Get-RGBValue 0xE7 0xE6 0xE6),  # lt2
    (Get-RGBValue 0x5B 0x9B 0xD5),  # accent1
    (Get-RGBValue 0xED 0x7D 0x31),  # accent2
    (Get-RGBValue 0xA5 0xA5 0xA5),  # accent3
    (Get-RGBValue 0xFF 0xC0 0x00),  # accent4
    (Get-RGBValue 0x44 0x72 0xC4),  # accent5
    (Get-RGBValue 0x70 0xAD 0x47),  # accent6
    (Get-RGBValue 0x05 0x63 0xC1),  # hlink
    (Get-RGBValue 0x95 0x4F 0x72)   # folHlink
)

$slide1 = $p.Slides.Item(1)
$themeColors = $slide1.ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = $officeColors[$i - 1]
}
